# Manejo de archivo de excel. Facturas enviadas a nueva hoja.
$wb = $excel.ActiveWorkbook

# Rename the original (only) sheet to "Pendientes"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Pendientes"

# Add a new sheet right after it, named "Facturados", and make it active
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Facturados"

# Copy the header row (with formatting) from Pendientes into Facturados
$ws1.Range("A1:V1").Copy($ws2.Range("A1"))

# Make the new sheet's header bold (matches the distinctive bold header style)
$ws2.Range("A1:V1").Font.Bold = $true

# Restore selection on each sheet to match the saved workbook view state
[void]$ws1.Range("P15").Select()
[void]$ws2.Range("R9").Select()

# Facturados is the active/visible tab when the workbook is reopened
[void]$ws2.Activate()
